$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two previously-empty cells with computed values
$ws.Range("C8").Value = 2710335.5680292002
$ws.Range("C11").Value = 2995998.4088369301

# Update the active selection to match the recorded cursor position
$ws.Range("C11").Select()
